$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Observed")

$ws.Range("E11").Value = 2.428060400516796
$ws.Range("E12").Value = 3.601082715295001
$ws.Range("E13").Value = 4.914803668672972
$ws.Range("E14").Value = 5.974550507929128
$ws.Range("E15").Value = 6.70534360788277
$ws.Range("E16").Value = 7.385609112332162
$ws.Range("E17").Value = 8.370199273368181
$ws.Range("E20").Value = 2.285
$ws.Range("E21").Value = 3.579664570230608
$ws.Range("E22").Value = 4.735967184801382
$ws.Range("E23").Value = 5.929332386363637
$ws.Range("E24").Value = 6.604126984126984
$ws.Range("E25").Value = 7.709947089947089
$ws.Range("E26").Value = 8.594074074074072
$ws.Range("E28").Value = 2.733696949168199
$ws.Range("E29").Value = 4.246589485395448
$ws.Range("E30").Value = 5.305275414013267
$ws.Range("E31").Value = 5.923018800622259
$ws.Range("E32").Value = 7.393218944221821
$ws.Range("E33").Value = 7.756830019760492
$ws.Range("E34").Value = 6.310178861788617
$ws.Range("E35").Value = 6.685714285714285
$ws.Range("E42").Value = 2.121700051647661
$ws.Range("E43").Value = 3.386942675159236
$ws.Range("E44").Value = 4.800179823737866
$ws.Range("E45").Value = 6.027890833527469
$ws.Range("E46").Value = 6.956493278341615
$ws.Range("E47").Value = 8.228275366779521
$ws.Range("E50").Value = 2.02566124813457
$ws.Range("E51").Value = 3.412988650693568
$ws.Range("E52").Value = 4.53639846743295
$ws.Range("E53").Value = 5.458722741433021
$ws.Range("E54").Value = 6.685076879747225
$ws.Range("E55").Value = 6.922885550245411
$ws.Range("E56").Value = 7.597444089456869
$ws.Range("E73").Value = 2.949916874480466
$ws.Range("E74").Value = 3.827930824100736
$ws.Range("E75").Value = 4.464450043933507
$ws.Range("E76").Value = 5.203025010647008
$ws.Range("E77").Value = 6.327372844651831
$ws.Range("E78").Value = 6.978159283685803
$ws.Range("E79").Value = 7.691452763876655
$ws.Range("E80").Value = 8.496567390283554
$ws.Range("E81").Value = 9.061967714048496
$ws.Range("E82").Value = 10.17996756690134
$ws.Range("E83").Value = 10.63377115675189
$ws.Range("E84").Value = 11.8807570434624
$ws.Range("E87").Value = 2.894669277264863
$ws.Range("E88").Value = 3.858043686174819
$ws.Range("E89").Value = 4.307022483286293
$ws.Range("E90").Value = 5.44758930969648
$ws.Range("E91").Value = 6.505628853397302
$ws.Range("E92").Value = 7.539993944227816
$ws.Range("E93").Value = 8.239207487680615
$ws.Range("E94").Value = 9.355048859934852
$ws.Range("E95").Value = 9.542792712173821
$ws.Range("E96").Value = 11.23461538461538
$ws.Range("E97").Value = 11.02202455524931
$ws.Range("E98").Value = 12.353125
$ws.Range("E101").Value = 3.077334005038731
$ws.Range("E102").Value = 4.171586641087303
$ws.Range("E103").Value = 4.941293109633294
$ws.Range("E104").Value = 5.662674900346253
$ws.Range("E105").Value = 6.955371699194276
$ws.Range("E106").Value = 7.621315192743762
$ws.Range("E107").Value = 8.555380972488388
$ws.Range("E108").Value = 9.246606282868557
$ws.Range("E109").Value = 9.569593147751606
$ws.Range("E110").ClearContents()
$ws.Range("E111").ClearContents()
$ws.Range("E112").Value = 13.75519480519481
$ws.Range("E113").Value = 14.05584415584416
$ws.Range("E114").Value = 15.97777777777778
$ws.Range("E119").Value = 3.580675740907558
$ws.Range("E120").Value = 3.649148345781931
$ws.Range("E121").Value = 4.834859223075668
$ws.Range("E122").Value = 6.568385122964684
$ws.Range("E123").Value = 8.289985371420514
$ws.Range("E132").Value = 2.79945611300716
$ws.Range("E133").Value = 3.881544929736314
$ws.Range("E134").Value = 4.518611548492419
$ws.Range("E135").Value = 5.183130897343855
$ws.Range("E136").Value = 6.315291904060707
$ws.Range("E137").Value = 6.78329595459256
$ws.Range("E138").Value = 7.493957516018458
$ws.Range("E139").Value = 8.34476246133107
$ws.Range("E140").Value = 8.905697177312417
$ws.Range("E141").Value = 10.34781949934124
$ws.Range("E142").Value = 10.11325428194993
$ws.Range("E153").Value = 2.77875
$ws.Range("E154").Value = 4.354672549019607
$ws.Range("E155").Value = 4.602333333333334
$ws.Range("E157").Value = 3.513480918489432
$ws.Range("E158").Value = 4.86938983845437
$ws.Range("E159").Value = 5.238481221538957
$ws.Range("E160").Value = 5.913443830570903
$ws.Range("E163").Value = 2.909907446501934
$ws.Range("E164").Value = 4.503084953055913
$ws.Range("E165").Value = 5.456542407603402
$ws.Range("E166").Value = 6.277260356720665
$ws.Range("E167").Value = 6.802236905136427
$ws.Range("E169").Value = 2.729770442922836
$ws.Range("E170").Value = 4.614949236065259
$ws.Range("E171").Value = 5.693032606646981
$ws.Range("E173").Value = 6.939375453885258
$ws.Range("E174").Value = 3.655773420479303
$ws.Range("E175").ClearContents()
$ws.Range("E176").ClearContents()
$ws.Range("E177").Value = 7.078125
$ws.Range("E192").Value = 1.005691339184544
$ws.Range("E193").Value = 2.844458049379705
$ws.Range("E194").Value = 3.875892214906213
$ws.Range("E195").Value = 5.194303974999472
$ws.Range("E196").Value = 6.779204073692867
$ws.Range("E198").Value = 1.532421694914059
$ws.Range("E199").Value = 2.616386858944526
$ws.Range("E200").Value = 4.972746375102701
$ws.Range("E201").Value = 4.899201261515694
$ws.Range("E202").Value = 6.058783771033148
$ws.Range("E203").Value = 7.106542821319646
$ws.Range("E204").Value = 8.479433878814683
$ws.Range("E205").Value = 9.34078584891982
